# Implementação do Extent Report
# - Update the leftover test-user string on Planilha1!A2
# - Switch the active sheet/selection from Planilha2!D2 to Planilha1!A2,
#   and move Planilha2's own lingering selection to D6.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Planilha1")
$ws2 = $wb.Worksheets.Item("Planilha2")

# Content change: the shared string "joao123998877" becomes "joao123998"
$ws1.Range("A2").Value = "joao123998"

# Move Planilha2's remembered selection (it is no longer the active tab)
$ws2.Range("D6").Select()

# Make Planilha1 the active/selected tab with A2 selected
$ws1.Activate()
$ws1.Range("A2").Select()
